$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5").Value = "152"
$ws.Range("D5").Value = "404886.40"
$ws.Range("C6").Value = "443"
$ws.Range("D6").Value = "1149510.82"
$ws.Range("C7").Value = "185"
$ws.Range("D7").Value = "411041.00"
$ws.Range("C8").Value = "809"
$ws.Range("D8").Value = "3074189.81"
$ws.Range("C9").Value = "32"
$ws.Range("D9").Value = "90600.00"
$ws.Range("C10").Value = "16"
$ws.Range("D10").Value = "41500.00"
$ws.Range("C13").Value = "99"
$ws.Range("D13").Value = "240800.00"
$ws.Range("C14").Value = "95"
$ws.Range("D14").Value = "231788.98"
$ws.Range("C16").Value = "136"
$ws.Range("D16").Value = "616217.26"
$ws.Range("C17").Value = "187"
$ws.Range("D17").Value = "423089.87"
$ws.Range("C23").Value = "314"
$ws.Range("D23").Value = "1293680.35"
$ws.Range("C37").Value = "379"
$ws.Range("D37").Value = "1511848.18"
$ws.Range("C74").Value = "21"
$ws.Range("D74").Value = "89400.00"
$ws.Range("C78").Value = "215"
$ws.Range("D78").Value = "595693.00"
$ws.Range("C79").Value = "26"
$ws.Range("D79").Value = "73991.00"
$ws.Range("C80").Value = "497"
$ws.Range("D80").Value = "2181324.03"
$ws.Range("C83").Value = "19"
$ws.Range("D83").Value = "60069.00"
$ws.Range("C85").Value = "37"
$ws.Range("D85").Value = "118669.00"
$ws.Range("C86").Value = "45"
$ws.Range("D86").Value = "103500.00"
$ws.Range("C89").Value = "108"
$ws.Range("D89").Value = "279620.00"
$ws.Range("C105").Value = "15"
$ws.Range("D105").Value = "37909.00"
$ws.Range("C106").Value = "23"
$ws.Range("D106").Value = "62209.84"
$ws.Range("C108").Value = "41"
$ws.Range("D108").Value = "128534.00"
$ws.Range("C109").Value = "19"
$ws.Range("D109").Value = "63913.61"
$ws.Range("C111").Value = "9"
$ws.Range("D111").Value = "34499.26"
$ws.Range("C112").Value = "6"
$ws.Range("D112").Value = "13500.00"
$ws.Range("C114").Value = "27"
$ws.Range("D114").Value = "73895.00"
$ws.Range("C115").Value = "16"
$ws.Range("D115").Value = "37100.00"
$ws.Range("C120").Value = "33"
$ws.Range("D120").Value = "130000.00"
$ws.Range("C121").Value = "65"
$ws.Range("D121").Value = "176877.00"
$ws.Range("C122").Value = "253"
$ws.Range("D122").Value = "706508.00"
$ws.Range("C123").Value = "123"
$ws.Range("D123").Value = "321012.45"
$ws.Range("C124").Value = "506"
$ws.Range("D124").Value = "2299836.06"
$ws.Range("C128").Value = "92"
$ws.Range("D128").Value = "278743.68"
$ws.Range("C133").Value = "129"
$ws.Range("D133").Value = "343995.68"
$ws.Range("C138").Value = "577"
$ws.Range("D138").Value = "1453546.00"
$ws.Range("C139").Value = "1842"
$ws.Range("D139").Value = "4965715.93"
$ws.Range("C140").Value = "2723"
$ws.Range("D140").Value = "6892397.55"
$ws.Range("C141").Value = "2698"
$ws.Range("D141").Value = "11810211.82"
$ws.Range("C142").Value = "365"
$ws.Range("D142").Value = "1047658.94"
$ws.Range("C144").Value = "258"
$ws.Range("D144").Value = "685516.33"
$ws.Range("C145").Value = "1065"
$ws.Range("D145").Value = "2803349.25"
$ws.Range("C146").Value = "519"
$ws.Range("D146").Value = "1553313.66"
$ws.Range("C147").Value = "383"
$ws.Range("D147").Value = "965377.83"
$ws.Range("C148").Value = "149"
$ws.Range("D148").Value = "369500.00"
$ws.Range("C149").Value = "434"
$ws.Range("D149").Value = "1421905.46"
$ws.Range("C150").Value = "864"
$ws.Range("D150").Value = "2104195.82"
$ws.Range("C197").Value = "56"
$ws.Range("D197").Value = "141153.50"
$ws.Range("C199").Value = "354"
$ws.Range("D199").Value = "953788.00"
